# eventbuttons.xlsx update
#  - HUB PWM "channel" terminology -> "port" (Commands sheet B66:C68)
#  - Commands sheet becomes the active / selected sheet (was Labels)
#  - Labels sheet scrolls back to the top (A1) and loses the tab-selected flag
#  - Commands sheet scrolls down near row 70 and the current selection moves to C66
#  - Commands!A column narrows from ~112.84 to ~26.74 characters

$wb = $excel.ActiveWorkbook

$labels   = $wb.Worksheets.Item("Labels")
$commands = $wb.Worksheets.Item("Commands")

# --- text fixes: channel -> port ------------------------------------------------
$commands.Cells.Item(66, 2).Value = "outhub(<port>,<value>[,<sn>])"
$commands.Cells.Item(66, 3).Value = "PHIDGET HUB PWM Output ON port <port> to  <value> in [0-100]"
$commands.Cells.Item(67, 2).Value = "togglehub(<port>[,<sn>])"
$commands.Cells.Item(67, 3).Value = "PHIDGET HUB PWM Output: toggles <port>"
$commands.Cells.Item(68, 2).Value = "pulsehub(<port>,<millis>[,<sn>])"
$commands.Cells.Item(68, 3).Value = "PHIDGET HUB PWM Output:  turn <port> ON for <millis> milliseconds"

# --- column width on Commands!A: 112.84 -> 26.74 --------------------------------
$commands.Columns.Item(1).ColumnWidth = 25.8

# --- view / selection state ------------------------------------------------------
# Labels was the active tab (topLeftCell A10, selection A28); it becomes inactive
# and scrolls back up to A1, keeping its A28 selection.
$labels.Activate()
$labels.Range("A28").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# Commands becomes the active tab, scrolled down (topLeftCell A70) with C66 selected.
$commands.Activate()
$commands.Range("C66").Select()
$excel.ActiveWindow.ScrollRow = 70
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "done"
